$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: "IS UNRESTRICTED" header ---
$ws.Range("L2").Value2 = "IS UNRESTRICTED"
$ws.Columns("L").ColumnWidth = 15.33203125

# --- Expand merged header / selection / view to include column L ---
$ws.Range("A1:K1").UnMerge()
$ws.Range("A1:L1").Merge()
$ws.Range("A1:L1").Select()
$ws.Application.ActiveWindow.Zoom = 120

# --- Row 1 fill (whole row, incl. future columns) + new row height ---
$ws.Rows(1).Interior.Color = 65535
$ws.Rows(1).RowHeight = 281

# --- Replace the big instructions cell (A1) with updated rich text ---
$full = @"
DO NOT DELETE THIS ROW! RETAIN THE HEADING ROW!
Instructions: Starting on Row 3, fill in the relevant fields. Do not delete rows 1 and 2.
For SPOUSE EMAIL and SPOUSE ID, these are mutually exclusive. If using them (they're not mandatory), use one or the other but not both. If both are used, the ID will take precendence and the email will be ignored.

NAME: Full Name
EMAIL: Propely formatted email address
MOBILE PHONE: In the format 04XXXXXXXX (spaces can be used)
GENDER: male, female, m or f
YEAR OF BIRTH: Optional. 4 digit year. Eg: 1985
APPOINTMENT: Optional. Allowed values only: elder, ministerial servant
SERVING AS: Optional. Allowed values only: field missionary, special pioneer, bethel family member, regular pioneer, publisher
MARITAL STATUS: Optional. Allowed values only: single, married, separated, divorced, widowed
SPOUSE EMAIL: Optional. Used to link spouses together. If a matching email is found, it will attach the users
SPOUSE ID: Optional. Used to link a user that already exists in the system to this user
RESPONSIBLE BROTHER: Inidcates in the system that a user (brother) has been trained to oversee a shift. Allowed values only. TRUE, FALSE.
IS UNRESTRICTED: TRUE is the default. If set to false (i.e. indicating they're a 'restricted' user), the volunteer cannot self-roster and they cannot see any shifts other than those they've been rostered onto. Allowed values only. TRUE, FALSE.
"@
$cell = $ws.Range("A1")
$cell.Value2 = $full

$boldRanges = @(@(1,352), @(354,5), @(370,6), @(409,13), @(469,7), @(498,14), @(546,12), @(617,11), @(744,15), @(837,13), @(947,10), @(1035,20), @(1173,16))
foreach ($r in $boldRanges) {
    $cell.Characters($r[0], $r[1]).Font.Bold = $true
}
